$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inventario")
$ws.Activate()

$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "3881"
$ws.Range("B8").Value = "Talco para bebé"
$ws.Range("C8").Value = 30
$ws.Range("D8").Value = 20

$ws.Range("D10").Select()
